$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3342.2778
$ws.Range("I74").Value = 3249.5
$ws.Range("J74").Value = 3458.25
$ws.Range("K74").Value = 3249.5
$ws.Range("L74").Value = 3458.25
$ws.Range("M74").Value = -2313.5
$ws.Range("N74").Value = -5330.25

$ws.Range("H77").Value = 3342.2778
$ws.Range("I77").Value = 3249.5
$ws.Range("J77").Value = 3458.25
$ws.Range("K77").Value = 16247.5
$ws.Range("L77").Value = 17291.25
$ws.Range("M77").Value = -11567.5
$ws.Range("N77").Value = -26651.25

$ws.Range("H132").Value = 4037561
$ws.Range("I132").Value = 4560.775
$ws.Range("J132").Value = 11370289
$ws.Range("K132").Value = 13682.325
$ws.Range("L132").Value = 34110867
$ws.Range("M132").Value = -11152.325
$ws.Range("N132").Value = -34115927

$ws.Range("H137").Value = 5265962.5
$ws.Range("I137").Value = 1519.875
$ws.Range("J137").Value = 14290722
$ws.Range("K137").Value = 4559.625
$ws.Range("L137").Value = 42872166
$ws.Range("M137").Value = -2009.625
$ws.Range("N137").Value = -42877266

$ws.Range("H138").Value = 3625183.8
$ws.Range("I138").Value = 1717.7838
$ws.Range("J138").Value = 7814816
$ws.Range("K138").Value = 5153.3514
$ws.Range("L138").Value = 23444448
$ws.Range("M138").Value = -13.35139999999956
$ws.Range("N138").Value = -23454728

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 5001801
$ws.Range("I132").Value = 7577060
$ws.Range("J132").Value = 2768.8823
$ws.Range("K132").Value = 22731180
$ws.Range("L132").Value = 8306.6469
$ws.Range("M132").Value = -22728650
$ws.Range("N132").Value = -13366.6469

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4066.394
$ws.Range("I134").Value = 3403.689
$ws.Range("J134").Value = 5486.476
$ws.Range("K134").Value = 10211.067
$ws.Range("L134").Value = 16459.428
$ws.Range("M134").Value = -7676.066999999999
$ws.Range("N134").Value = -21529.428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2559.0715
$ws.Range("I132").Value = 1988.174
$ws.Range("J132").Value = 5185.2
$ws.Range("K132").Value = 5964.522
$ws.Range("L132").Value = 15555.6
$ws.Range("M132").Value = -3434.522
$ws.Range("N132").Value = -20615.6

$ws.Range("H134").Value = 487575.94
$ws.Range("I134").Value = 1493.3823
$ws.Range("J134").Value = 1589363.1
$ws.Range("K134").Value = 4480.1469
$ws.Range("L134").Value = 4768089.300000001
$ws.Range("M134").Value = -1945.1469
$ws.Range("N134").Value = -4773159.300000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 6165.3887
$ws.Range("I3").Value = 4593.125
$ws.Range("J3").Value = 7423.2
$ws.Range("K3").Value = 13779.375
$ws.Range("L3").Value = 22269.6
$ws.Range("M3").Value = -13667.375
$ws.Range("N3").Value = -22493.6

$ws.Range("H5").Value = 1527.3077
$ws.Range("I5").Value = 562.61536
$ws.Range("J5").Value = 2492
$ws.Range("K5").Value = 1687.84608
$ws.Range("L5").Value = 7476
$ws.Range("M5").Value = -1575.84608
$ws.Range("N5").Value = -7700

$ws.Range("H33").Value = 721589.6
$ws.Range("I33").Value = 1122412.5
$ws.Range("J33").Value = 108.4
$ws.Range("K33").Value = 6734475
$ws.Range("L33").Value = 650.4000000000001
$ws.Range("M33").Value = -6734192
$ws.Range("N33").Value = -1216.4

$ws.Range("H34").Value = 1180.1578
$ws.Range("J34").Value = 1229.0555
$ws.Range("L34").Value = 3687.1665
$ws.Range("N34").Value = -3855.1665

$ws.Range("H38").Value = 161.5
$ws.Range("J38").Value = 54.666668
$ws.Range("L38").Value = 164.000004
$ws.Range("N38").Value = -858.000004

$ws.Range("H39").Value = 504.9524
$ws.Range("J39").Value = 504.9524
$ws.Range("L39").Value = 1514.8572
$ws.Range("N39").Value = -2102.8572

$ws.Range("H40").Value = 145.36363
$ws.Range("I40").Value = 112.5
$ws.Range("J40").Value = 233
$ws.Range("K40").Value = 450
$ws.Range("L40").Value = 932
$ws.Range("M40").Value = -381
$ws.Range("N40").Value = -1070

$ws.Range("H109").Value = 2840.7083
$ws.Range("I109").Value = 1575.2222
$ws.Range("J109").Value = 3600
$ws.Range("K109").Value = 4725.6666
$ws.Range("L109").Value = 10800
$ws.Range("M109").Value = -3685.6666
$ws.Range("N109").Value = -12880

$ws.Range("H125").Value = 4366.6
$ws.Range("J125").Value = 5208.25
$ws.Range("L125").Value = 15624.75
$ws.Range("N125").Value = -25464.75

$ws.Range("H132").Value = 2299.7856
$ws.Range("I132").Value = 1000
$ws.Range("J132").Value = 3274.625
$ws.Range("K132").Value = 9000
$ws.Range("L132").Value = 29471.625
$ws.Range("M132").Value = -6470
$ws.Range("N132").Value = -34531.625

$ws.Range("H133").Value = 153850610
$ws.Range("I133").Value = 222225360
$ws.Range("J133").Value = 7444
$ws.Range("K133").Value = 666676080
$ws.Range("L133").Value = 22332
$ws.Range("M133").Value = -666671020
$ws.Range("N133").Value = -32452

$ws.Range("H134").Value = 4406.25
$ws.Range("I134").Value = 2839.1667
$ws.Range("J134").Value = 5973.3335
$ws.Range("K134").Value = 8517.500100000001
$ws.Range("L134").Value = 17920.0005
$ws.Range("M134").Value = -3447.500100000001
$ws.Range("N134").Value = -28060.0005

$ws.Range("H135").Value = 1527.3077
$ws.Range("I135").Value = 562.61536
$ws.Range("J135").Value = 2492
$ws.Range("K135").Value = 5063.53824
$ws.Range("L135").Value = 22428
$ws.Range("M135").Value = -2528.53824
$ws.Range("N135").Value = -27498

$ws.Range("H136").Value = 4804.2104
$ws.Range("I136").Value = 3195
$ws.Range("J136").Value = 5233.3335
$ws.Range("K136").Value = 9585
$ws.Range("L136").Value = 15700.0005
$ws.Range("M136").Value = -4485
$ws.Range("N136").Value = -25900.0005

$ws.Range("H138").Value = 7101
$ws.Range("I138").Value = 2977.5
$ws.Range("J138").Value = 9850
$ws.Range("K138").Value = 8932.5
$ws.Range("L138").Value = 29550
$ws.Range("M138").Value = -3792.5
$ws.Range("N138").Value = -39830

$ws.Range("H139").Value = 2021.3103
$ws.Range("I139").Value = 1414
$ws.Range("J139").Value = 2450
$ws.Range("K139").Value = 4242
$ws.Range("L139").Value = 7350
$ws.Range("M139").Value = 898
$ws.Range("N139").Value = -17630

$ws.Range("H140").Value = 3488.484
$ws.Range("I140").Value = 1502.3158
$ws.Range("K140").Value = 4506.9474
$ws.Range("M140").Value = 673.0526

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3783.4324
$ws.Range("I132").Value = 2435.9167
$ws.Range("J132").Value = 6271.154
$ws.Range("K132").Value = 7307.750100000001
$ws.Range("L132").Value = 18813.462
$ws.Range("M132").Value = -4777.750100000001
$ws.Range("N132").Value = -23873.462

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5395.758
$ws.Range("I122").Value = 6233.6665
$ws.Range("J122").Value = 4390.2666
$ws.Range("K122").Value = 18700.9995
$ws.Range("L122").Value = 13170.7998
$ws.Range("M122").Value = -16250.9995
$ws.Range("N122").Value = -18070.7998

$ws.Range("H132").Value = 7818351.5
$ws.Range("I132").Value = 3665.3953
$ws.Range("J132").Value = 23819852
$ws.Range("K132").Value = 10996.1859
$ws.Range("L132").Value = 71459556
$ws.Range("M132").Value = -8466.1859
$ws.Range("N132").Value = -71464616

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 11557.357
$ws.Range("J62").Value = 11984.846
$ws.Range("L62").Value = 11984.846
$ws.Range("N62").Value = -13232.846

$ws.Range("H65").Value = 11557.357
$ws.Range("J65").Value = 11984.846
$ws.Range("L65").Value = 59924.23
$ws.Range("N65").Value = -66164.23

$ws.Range("H132").Value = 6197.926
$ws.Range("I132").Value = 7438.9473
$ws.Range("J132").Value = 3250.5
$ws.Range("K132").Value = 22316.8419
$ws.Range("L132").Value = 9751.5
$ws.Range("M132").Value = -19786.8419
$ws.Range("N132").Value = -14811.5

